$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (Excel serial 45206 = 2023-10-07)
# for every data row (2 through 473). The edit shifts that date by 2 days to
# serial 45208 (2023-10-09) for all of those rows.
$ws.Range("C2:C473").Value = 45208
